$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.823.02'
$ws.Range("E2").Value = '  +0.29%  '

$ws.Range("D3").Value = '2.565.77'
$ws.Range("E3").Value = '  +1.27%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.32'
$ws.Range("E5").Value = '  -1.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.46'
$ws.Range("E6").Value = '  +3.07%  '

$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.63'
$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("E12").Value = '  -1.00%  '

$ws.Range("D13").Value = '2.961.03'
$ws.Range("E13").Value = '  +1.41%  '

$ws.Range("E14").Value = '  -1.50%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.83'
$ws.Range("E15").Value = '  +5.32%  '

$ws.Range("D16").Value = '2.545.01'
$ws.Range("E16").Value = '  +1.68%  '

$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '42.852.31'
$ws.Range("E18").Value = '  +0.19%  '

$ws.Range("D20").Value = '0.0₃0961'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.43'
$ws.Range("E21").Value = '  -2.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '69.44'
$ws.Range("E22").Value = '  -0.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '248.41'
$ws.Range("E23").Value = '  -0.96%  '

$ws.Range("E24").Value = '  -0.12%  '

$ws.Range("E25").Value = '  +0.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '27.02'
$ws.Range("E26").Value = '  +2.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.39'
$ws.Range("E28").Value = '  -0.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '39.94'
$ws.Range("E29").Value = '  -0.81%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '10.17'
$ws.Range("E30").Value = '  -1.93%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '159.18'
$ws.Range("E31").Value = '  +1.79%  '

$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0797'
$ws.Range("E33").Value = '  +2.42%  '

$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.10'
$ws.Range("E34").Value = '  -2.14%  '

$ws.Range("B35").Value = 'WEMIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.67'
$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("E36").Value = '  -1.15%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.61'
$ws.Range("E37").Value = '  -0.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.57'

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.111'
$ws.Range("E39").Value = '  +0.84%  '

$ws.Range("E40").Value = '  -0.21%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '22.87'
$ws.Range("E41").Value = '  +2.16%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.09'
$ws.Range("E42").Value = '  +7.34%  '

$ws.Range("E43").Value = '  -0.12%  '

$ws.Range("E44").Value = '  -0.37%  '

$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.21'
$ws.Range("E45").Value = '  -1.03%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.993.87'
$ws.Range("E46").Value = '  -1.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.01'
$ws.Range("E47").Value = '  -0.51%  '

$ws.Range("D48").Value = '2.812.05'
$ws.Range("E48").Value = '  +1.44%  '

$ws.Range("E49").Value = '  +2.82%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '81.35'
$ws.Range("E50").Value = '  -3.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.05'
$ws.Range("E51").Value = '  -1.32%  '
